$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.08912883471871574
$ws.Range("E2").Value = 0.01050834099566544

$ws.Range("D3").Value = 0.1057425340677712
$ws.Range("E3").Value = 0.002067749194186064

$ws.Range("D4").Value = 0.1206448017561099
$ws.Range("E4").Value = -0.008569647865378616

$ws.Range("D5").Value = 0.1419200458636643
$ws.Range("E5").Value = -0.008848858339258858

$ws.Range("D6").Value = 0.1384851453024938
$ws.Range("E6").Value = -0.009722032041626916

$ws.Range("D7").Value = 0.1483726659375793
$ws.Range("E7").Value = -0.01297610156833462

$ws.Range("D8").Value = 0.1267092620719243
$ws.Range("E8").Value = -0.009877282250823072

$ws.Range("D9").Value = 0.1289967102817416
$ws.Range("E9").Value = -0.0009389976707979697

$ws.Range("E10").Value = -0.005778795180670371

$ws.Protect()
